$d = $word.ActiveDocument

# The document starts with three "Title"-styled paragraphs:
#   1. empty paragraph carrying bookmark "_81e444s2jpyq"
#   2. empty paragraph carrying bookmark "_iy39xx24l9dl"
#   3. paragraph with text "Declaração do Problema" carrying bookmark "_htklu03vixnz"
# The edit removes the two leading empty Title paragraphs, leaving the
# "Declaração do Problema" paragraph (with its bookmark) as the first
# paragraph of the document.
$d.Paragraphs(2).Range.Delete()
$d.Paragraphs(1).Range.Delete()

# Merge the three runs "é o prejuízo causado pela" + " " + "alta demanda de
# pedidos que " into a single replacement run reading
# "da falta de administração e organização de pedidos ".
$d.Content.Find.Execute(
    "é o prejuízo causado pela alta demanda de pedidos que ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "da falta de administração e organização de pedidos ", 2)

# Replace the closing sentence fragment.
$d.Content.Find.Execute(
    "a falta de controle financeiro e perda de clientes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a perda de clientes e lucro para a empresa.", 2)
